# "reg 302 e 303 terminadas"
# Adds two new registry rows (REG-402 / REG-403) to Sheet1, plus a new
# "in_Quantidade_Parcelas" column (I), and two stray underlined cells
# (H9, G10) left over from formatting further down the sheet.
#
# NOTE: the order in which new string values are first assigned controls
# the order they are appended to xl/sharedStrings.xml, so the assignment
# order below intentionally matches the order of the new <si> entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: REG-402 block
$ws.Range("A5").Value = "REG-402"
$ws.Range("A5").Font.Underline = 2          # xlUnderlineStyleSingle
$ws.Range("H5").Value = "'3250000"          # stored as text (quote-prefixed number)

# New column header (I1) + its row-5 value
$ws.Range("I1").Value = "in_Quantidade_Parcelas"
$ws.Range("I5").Value = "'5"                # stored as text (quote-prefixed number)

# Row 6: REG-403 block
$ws.Range("A6").Value = "REG-403"
$ws.Range("H6").Value = "'6963455"          # stored as text (quote-prefixed number)
$ws.Range("H6").Font.Underline = 2          # xlUnderlineStyleSingle
$ws.Range("I6").Value = "'2"                # stored as text (quote-prefixed number)

# Two empty, underlined-format-only cells further down the sheet
$ws.Range("H9").Font.Underline = 2
$ws.Range("G10").Font.Underline = 2

# Leave the selection where the author left it when they saved
$ws.Range("I11").Select() | Out-Null
